# Add contributor info (name, email, repo link) to Sheet1, row 2,
# matching the header row ("name", "email", "Repo Link") already in A1:C1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 values (write right-to-left so the shared-string table ends up
#     ordered repo-link, email, name - matching the source workbook) -------
$name     = "السيد اسامه رجب السيد"
$email    = "heikalsayed@gmail.com"
$repoUrl  = "https://github.com/0xkillua/Security-Task.git"

$ws.Range("C2").Value = $repoUrl
$ws.Range("B2").Value = $email
$ws.Range("A2").Value = $name

# --- Hyperlinks (repo link first, then email, to mirror rId ordering) -----
$ws.Hyperlinks.Add($ws.Range("C2"), $repoUrl) | Out-Null
$ws.Hyperlinks.Add($ws.Range("B2"), "mailto:" + $email) | Out-Null

# --- Widen the columns so the new long text fits (values approximate the
#     auto-fit widths recorded in the target file) -------------------------
$ws.Columns.Item(1).ColumnWidth = 45.6
$ws.Columns.Item(2).ColumnWidth = 48.8
$ws.Columns.Item(3).ColumnWidth = 50.8

# --- Put the selection on A2, like in the saved workbook ------------------
$ws.Range("A2").Select() | Out-Null
